$d = $word.ActiveDocument

# Sorcière section: the two sentences describing the death/life potions were
# split across two runs ("... et " + "pourcentage ..."). Word's Find/Replace
# merges any text it matches (even across run boundaries) into a single run,
# so searching for the already-concatenated text and "replacing" it with the
# same text forces the two runs to collapse into one, matching the target
# OOXML (single <w:r> per bullet).
$d.Content.Find.Execute("nombre de fois la potion de mort a été utilisé et pourcentage de loup-garou tuer", $true, $false, $false, $false, $false, $true, 1, $false, "nombre de fois la potion de mort a été utilisé et pourcentage de loup-garou tuer", 2)

$d.Content.Find.Execute("nombre de fois la potion de vie a été utilisé et pourcentage où elle s’est sauvé et ou elle a sauvé quelqu’un d’autre", $true, $false, $false, $false, $false, $true, 1, $false, "nombre de fois la potion de vie a été utilisé et pourcentage où elle s’est sauvé et ou elle a sauvé quelqu’un d’autre", 2)

# Salvateur section: drop the "lg protégée et de" part, keeping only "villageois".
$d.Content.Find.Execute("Pourcentage de lg protégée et de villageois", $true, $false, $false, $false, $false, $true, 1, $false, "Pourcentage de villageois", 2)
